$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last status check" timestamp in the header (F1).
$ws.Range("F1").Value = "Last status check on: 01.02.2022 06:00"

# Row 4 (Globus): new scrape result.
# - B4 (Cena / current price): new reading 35.7
# - C4 (Old Cena): previous reading (what used to be in B4)
# - D4 (Delta Cena): written by the scraper as signed text "+0.4"
# - E4 (Old Datum): written by the scraper as a plain timestamp string
#   instead of a serial date, and loses the date number-format style.

# Shift the previous "current price" (B4) into the "old price" slot (C4)
# before overwriting B4 with the freshly scraped price. Use Value2 for the
# read - Value's getter is unreliable in this bridge.
$ws.Range("C4").Value2 = $ws.Range("B4").Value2
$ws.Range("B4").Value2 = 35.7

# D4 must be literal text "+0.4", not a recalculated number - force text
# formatting before the write so the "+" sign survives, then drop back to
# the workbook's default (unstyled) cell style to match the target file.
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "+0.4"
$ws.Range("D4").Style = "Normal"

# E4 becomes a plain text timestamp (no more date serial / date style).
$ws.Range("E4").Value = "2022-02-01 06:00:13"
$ws.Range("E4").Style = "Normal"
